# Updated all assays to accommodate the newly introduced dataset type.
#
# The "dataset_type" lookup sheet gets:
#   - "nanoPOTS" row removed
#   - "NanoDESI" row removed
#   - a new "2D Imaging Mass Cytometry" row added right after "MALDI"
# and the Visium sheet's dataset_type validation range shrinks to match
# the (now 35-row) lookup list. The template's createdOn timestamp is
# also bumped to record the edit.

$wb = $excel.ActiveWorkbook

$dsType = $wb.Worksheets.Item("dataset_type")

# Remove "nanoPOTS" (row 3).
$dsType.Rows.Item(3).Delete()

# Remove "NanoDESI" (was row 21, now row 20 after the delete above).
$dsType.Rows.Item(20).Delete()

# "MALDI" is now row 22; insert the new assay right after it.
$dsType.Rows.Item(23).Insert()
$dsType.Cells.Item(23, 1).Value = "2D Imaging Mass Cytometry"
$dsType.Cells.Item(23, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000296"

# The lookup list is now 35 rows (was 36); shrink the Visium sheet's
# dataset_type validation range to match.
$visium = $wb.Worksheets.Item("Visium")
$visium.Range("D2:D1001").Validation.Formula1 = "'dataset_type'!`$A`$1:`$A`$35"

# Bump the template's recorded edit timestamp.
$metadata = $wb.Worksheets.Item(".metadata")
$metadata.Cells.Item(2, 3).Value = "2023-11-02T15:46:14-07:00"
